$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new reference-sequence / name lookup rows to the "sequences"
#    sheet (rows 47-49), introducing the new alternative postfix sequences.
# ---------------------------------------------------------------------------
$seq = $wb.Worksheets.Item("sequences")

$seq.Range("A47").Value = "i57_chunk"
$seq.Range("B47").Value = "GCAGGGCGGTTTTTCGAAGGTTCTCTGAGCTACCAACTCTTTGAACCG"

$seq.Range("A48").Value = "u21attb"
$seq.Range("B48").Value = "TCCGTCTACGAACTCCCAGCAGGTAGGTATGATCCTGACGACGGAGCACGCCGTCGTCGACAAGCC"

$seq.Range("A49").Value = "u22attb"
$seq.Range("B49").Value = "GCTTGGATTCTGCGTTTGTTAGGTATGATCCTGACGACGGAGCACGCCGTCGTCGACAAGCC"

# ---------------------------------------------------------------------------
# 2. Update the "alldata" sheet: column N (prefix) and column Q (suffix)
#    for rows 10-151 now reference the new, more specific sequence names
#    instead of the generic "attBnoU" / "genomechunk" placeholders.
# ---------------------------------------------------------------------------
$alldata = $wb.Worksheets.Item("alldata")

$rowUpdates = @(
    "10|u21attb|",
    "11|u21attb|",
    "12|u21attb|",
    "13|u21attb|",
    "14|u21attb|",
    "15|u22attb|",
    "16|u22attb|",
    "17|u22attb|",
    "18|u22attb|",
    "19|u22attb|",
    "20|u22attb|",
    "21|u22attb|",
    "22|u22attb|",
    "23|u22attb|",
    "24|u22attb|",
    "25|u22attb|",
    "26|u22attb|",
    "27|u22attb|",
    "28|u22attb|",
    "29|u22attb|",
    "30|u22attb|",
    "31|u22attb|",
    "32|u22attb|",
    "33|u22attb|",
    "34|u22attb|",
    "35|u22attb|",
    "36|u22attb|",
    "37|u22attb|",
    "38|u22attb|",
    "39|u22attb|",
    "40|u22attb|",
    "41|u22attb|",
    "42|u21attb|",
    "43|u21attb|",
    "44|u21attb|",
    "45|u21attb|",
    "46|u21attb|",
    "47|u21attb|",
    "48|u21attb|",
    "49|u21attb|",
    "50|u21attb|",
    "51|u21attb|",
    "52|u21attb|",
    "53|u21attb|",
    "54|u21attb|",
    "55|u21attb|",
    "56|u21attb|",
    "57|u21attb|i57_chunk",
    "58|u21attb|i57_chunk",
    "59|u21attb|i57_chunk",
    "60|u21attb|i57_chunk",
    "61|u21attb|i57_chunk",
    "62|u21attb|i57_chunk",
    "63|u21attb|i57_chunk",
    "64|u21attb|i57_chunk",
    "65|u21attb|i57_chunk",
    "66|u21attb|",
    "67|u21attb|",
    "68|u21attb|",
    "69|u21attb|",
    "70|u21attb|",
    "71|u21attb|",
    "72|u21attb|",
    "73|u21attb|",
    "74|u21attb|",
    "75|u21attb|[genomechunk,i57_chunk]",
    "76|u21attb|[genomechunk,i57_chunk]",
    "77|u21attb|[genomechunk,i57_chunk]",
    "78|u21attb|[genomechunk,i57_chunk]",
    "79|u21attb|[genomechunk,i57_chunk]",
    "80|u21attb|[genomechunk,i57_chunk]",
    "81|u21attb|[genomechunk,i57_chunk]",
    "82|u21attb|[genomechunk,i57_chunk]",
    "83|u21attb|[genomechunk,i57_chunk]",
    "84|u21attb|[genomechunk,i57_chunk]",
    "85|u21attb|[genomechunk,i57_chunk]",
    "86|u21attb|[genomechunk,i57_chunk]",
    "87|u21attb|[genomechunk,i57_chunk]",
    "88|u21attb|[genomechunk,i57_chunk]",
    "89|u21attb|[genomechunk,i57_chunk]",
    "90|u21attb|[genomechunk,i57_chunk]",
    "91|u21attb|[genomechunk,i57_chunk]",
    "92|u21attb|[genomechunk,i57_chunk]",
    "93|u21attb|[genomechunk,i57_chunk]",
    "94|u21attb|[genomechunk,i57_chunk]",
    "95|u21attb|[genomechunk,i57_chunk]",
    "96|u21attb|[genomechunk,i57_chunk]",
    "97|u21attb|[genomechunk,i57_chunk]",
    "98|u21attb|[genomechunk,i57_chunk]",
    "99|u21attb|[genomechunk,i57_chunk]",
    "100|u21attb|[genomechunk,i57_chunk]",
    "101|u21attb|[genomechunk,i57_chunk]",
    "102|u21attb|[genomechunk,i57_chunk]",
    "103|u21attb|[genomechunk,i57_chunk]",
    "104|u21attb|[genomechunk,i57_chunk]",
    "105|u21attb|[genomechunk,i57_chunk]",
    "106|u21attb|[genomechunk,i57_chunk]",
    "107|u21attb|[genomechunk,i57_chunk]",
    "108|u21attb|[genomechunk,i57_chunk]",
    "109|u21attb|[genomechunk,i57_chunk]",
    "110|u21attb|[genomechunk,i57_chunk]",
    "111|u21attb|[genomechunk,i57_chunk]",
    "112|u21attb|[genomechunk,i57_chunk]",
    "113|u21attb|[genomechunk,i57_chunk]",
    "114|u21attb|[genomechunk,i57_chunk]",
    "115|u21attb|[genomechunk,i57_chunk]",
    "116|u21attb|[genomechunk,i57_chunk]",
    "117|u21attb|[genomechunk,i57_chunk]",
    "118|u21attb|[genomechunk,i57_chunk]",
    "119|u21attb|[genomechunk,i57_chunk]",
    "120|u21attb|[genomechunk,i57_chunk]",
    "121|u21attb|[genomechunk,i57_chunk]",
    "122|u21attb|[genomechunk,i57_chunk]",
    "123|u21attb|[genomechunk,i57_chunk]",
    "124|u21attb|[genomechunk,i57_chunk]",
    "125|u21attb|[genomechunk,i57_chunk]",
    "126|u21attb|[genomechunk,i57_chunk]",
    "127|u21attb|[genomechunk,i57_chunk]",
    "128|u22attb|[genomechunk,i57_chunk]",
    "129|u21attb|[genomechunk,i57_chunk]",
    "130|u21attb|[genomechunk,i57_chunk]",
    "131|u22attb|[genomechunk,i57_chunk]",
    "132|u21attb|[genomechunk,i57_chunk]",
    "133|u21attb|[genomechunk,i57_chunk]",
    "134|u22attb|[genomechunk,i57_chunk]",
    "135|u21attb|[genomechunk,i57_chunk]",
    "136|u21attb|[genomechunk,i57_chunk]",
    "137|u22attb|[genomechunk,i57_chunk]",
    "138|u21attb|[genomechunk,i57_chunk]",
    "139|u21attb|[genomechunk,i57_chunk]",
    "140|u22attb|",
    "141|u21attb|",
    "142|u21attb|",
    "143|u22attb|",
    "144|u21attb|",
    "145|u21attb|",
    "146|u22attb|",
    "147|u21attb|",
    "148|u21attb|",
    "149|u22attb|",
    "150|u21attb|",
    "151|u21attb|"
)

foreach ($entry in $rowUpdates) {
    $parts = $entry.Split("|")
    $r = [int]$parts[0]
    $nVal = $parts[1]
    $qVal = $parts[2]

    if ($nVal -ne "") {
        $alldata.Cells.Item($r, 14).Value = $nVal
    }
    if ($qVal -ne "") {
        $alldata.Cells.Item($r, 17).Value = $qVal
    }
}

# ---------------------------------------------------------------------------
# 3. Restore the sheet-view scroll/selection state recorded in the workbook
#    (frozen-pane top-left cell and active selection) for both sheets.
# ---------------------------------------------------------------------------
$alldataView = $alldata.Application.ActiveWindow
$alldata.Activate()
$excel.ActiveWindow.ScrollRow = 107
$alldata.Range("W129").Select()

$seq.Activate()
$excel.ActiveWindow.ScrollRow = 28
$seq.Range("F50").Select()

$alldata.Activate()
